# Corpus/Default_Messages.xlsx — update the welcome message text and let
# the row re-fit to its (now shorter) two-line wrapped content.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B2 holds the SAM welcome message (A2 = "welcome" / Message_Type).
$ws.Range("B2").Value = "Hi, I am SAM, the Search Assistant Manager. How can I help you today ?"

# The cell still wraps (style carried over unchanged); the shorter text
# only needs two lines instead of three, so the row shrinks accordingly.
$ws.Rows.Item(2).RowHeight = 28.8
